$wb = $excel.ActiveWorkbook

$successMsg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet "Bico" ---
$wsBico = $wb.Worksheets.Item("Bico")

$wsBico.Range("H1").Value = "Obs_relatorio"
$wsBico.Range("I1").Value = "Obs_sped"

for ($r = 2; $r -le 13; $r++) {
    $wsBico.Cells.Item($r, 8).Value = $successMsg
    $wsBico.Cells.Item($r, 9).Value = "'"
    $wsBico.Cells.Item($r, 9).Style = "Normal"
}

# --- Sheet "Tanque" ---
$wsTanque = $wb.Worksheets.Item("Tanque")

$wsTanque.Range("F1").Value = "Obs_relatorio"
$wsTanque.Range("G1").Value = "Obs_sped"

for ($r = 2; $r -le 8; $r++) {
    $wsTanque.Cells.Item($r, 6).Value = $successMsg
    $wsTanque.Cells.Item($r, 7).Value = "'"
    $wsTanque.Cells.Item($r, 7).Style = "Normal"
}
